# Insert a new data row at row 414 (pushing existing rows 414-469 down to 415-470)
# and populate it with the new record described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 414; this shifts rows 414..469 down to 415..470
# and copies formatting (e.g. the date style on column D) from the row being pushed down.
$ws.Rows.Item(414).Insert()

# Populate the newly inserted row 414 with the new record's data.
$ws.Range("A414").Value = 5
$ws.Range("B414").Value = "Macroferia Regional de Talca"
$ws.Range("C414").Value = "Maule"
$ws.Range("D414").Value = 45142
$ws.Range("E414").Value = 7
$ws.Range("F414").Value = 100112045
$ws.Range("G414").Value = "Zapallo"
$ws.Range("H414").Value = "Camote"
$ws.Range("I414").Value = "1a (guarda)"
$ws.Range("J414").Value = 800
$ws.Range("K414").Value = 350
$ws.Range("L414").Value = 350
$ws.Range("M414").Value = 350
$ws.Range("N414").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O414").Value = "Región del Maule"
$ws.Range("P414").Value = 350
$ws.Range("Q414").Value = 1
$ws.Range("R414").Value = "Hortaliza"
